$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.84976866666667
$ws.Range("H2").Value = 65.549306
$ws.Range("I2").Value = 0.05020018890879543
$ws.Range("J2").Value = 0.05020018890879543
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.231909
$ws.Range("N2").Value = 24.695727
$ws.Range("O2").Value = 0.6819745823584403
$ws.Range("P2").Value = 0.6819745823584403
$ws.Range("Q2").Value = 179.8653073350513
$ws.Range("R2").Value = 1618.787766015462
$ws.Range("S2").Value = 0.03423525286539057
$ws.Range("T2").Value = 0.03423525286539057
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.84976866666667
$ws.Range("H3").Value = 65.549306
$ws.Range("I3").Value = 0.05020018890879543
$ws.Range("J3").Value = 0.05020018890879543
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.472366666666667
$ws.Range("N3").Value = 4.4171
$ws.Range("O3").Value = 0.1219785887548671
$ws.Range("P3").Value = 0.121978588754867
$ws.Range("Q3").Value = 32.17087105917778
$ws.Range("R3").Value = 289.5378395326
$ws.Range("S3").Value = 0.006123348198322596
$ws.Range("T3").Value = 0.006123348198322594
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.84976866666667
$ws.Range("H4").Value = 65.549306
$ws.Range("I4").Value = 0.05020018890879543
$ws.Range("J4").Value = 0.05020018890879543
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.366422
$ws.Range("N4").Value = 7.099266000000001
$ws.Range("O4").Value = 0.1960468288866927
$ws.Range("P4").Value = 0.1960468288866926
$ws.Range("Q4").Value = 51.70577326771068
$ws.Range("R4").Value = 465.351959409396
$ws.Range("S4").Value = 0.009841587845082265
$ws.Range("T4").Value = 0.009841587845082263
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 385.0524703333334
$ws.Range("H5").Value = 1155.157411
$ws.Range("I5").Value = 0.8846641374295412
$ws.Range("J5").Value = 0.8846641374295412
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.231909
$ws.Range("N5").Value = 24.695727
$ws.Range("O5").Value = 0.6819745823584403
$ws.Range("P5").Value = 0.6819745823584403
$ws.Range("Q5").Value = 3169.7168960092
$ws.Range("R5").Value = 28527.4520640828
$ws.Range("S5").Value = 0.6033184556510012
$ws.Range("T5").Value = 0.6033184556510012
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 385.0524703333334
$ws.Range("H6").Value = 1155.157411
$ws.Range("I6").Value = 0.8846641374295412
$ws.Range("J6").Value = 0.8846641374295412
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.472366666666667
$ws.Range("N6").Value = 4.4171
$ws.Range("O6").Value = 0.1219785887548671
$ws.Range("P6").Value = 0.121978588754867
$ws.Range("Q6").Value = 566.9384222364556
$ws.Range("R6").Value = 5102.4458001281
$ws.Range("S6").Value = 0.1079100830056972
$ws.Range("T6").Value = 0.1079100830056972
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 385.0524703333334
$ws.Range("H7").Value = 1155.157411
$ws.Range("I7").Value = 0.8846641374295412
$ws.Range("J7").Value = 0.8846641374295412
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.366422
$ws.Range("N7").Value = 7.099266000000001
$ws.Range("O7").Value = 0.1960468288866927
$ws.Range("P7").Value = 0.1960468288866926
$ws.Range("Q7").Value = 911.1966369511475
$ws.Range("R7").Value = 8200.769732560328
$ws.Range("S7").Value = 0.1734355987728428
$ws.Range("T7").Value = 0.1734355987728428
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 28.350479
$ws.Range("H8").Value = 85.05143699999999
$ws.Range("I8").Value = 0.06513567366166337
$ws.Range("J8").Value = 0.06513567366166337
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 8.231909
$ws.Range("N8").Value = 24.695727
$ws.Range("O8").Value = 0.6819745823584403
$ws.Range("P8").Value = 0.6819745823584403
$ws.Range("Q8").Value = 233.378563234411
$ws.Range("R8").Value = 2100.407069109699
$ws.Range("S8").Value = 0.04442087384204853
$ws.Range("T8").Value = 0.04442087384204853
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 28.350479
$ws.Range("H9").Value = 85.05143699999999
$ws.Range("I9").Value = 0.06513567366166337
$ws.Range("J9").Value = 0.06513567366166337
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.472366666666667
$ws.Range("N9").Value = 4.4171
$ws.Range("O9").Value = 0.1219785887548671
$ws.Range("P9").Value = 0.121978588754867
$ws.Range("Q9").Value = 41.74230026363333
$ws.Range("R9").Value = 375.6807023727
$ws.Range("S9").Value = 0.007945157550847261
$ws.Range("T9").Value = 0.00794515755084726
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 28.350479
$ws.Range("H10").Value = 85.05143699999999
$ws.Range("I10").Value = 0.06513567366166337
$ws.Range("J10").Value = 0.06513567366166337
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.366422
$ws.Range("N10").Value = 7.099266000000001
$ws.Range("O10").Value = 0.1960468288866927
$ws.Range("P10").Value = 0.1960468288866926
$ws.Range("Q10").Value = 67.089197216138
$ws.Range("R10").Value = 603.802774945242
$ws.Range("S10").Value = 0.01276964226876757
$ws.Range("T10").Value = 0.01276964226876757

Write-Output "Applied NATMI value updates"